# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new labels in AC1:AE1, formatted like the other headers.
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header formatting (bold font, border, centered alignment) from an
# existing header cell onto the new header cells without touching their values.
$headerFormat = $ws.Range("A1")
$headerFormat.Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Data rows (2-46): every team had the same 1990 season record.
$ws.Range("AC2:AC46").Value = 91
$ws.Range("AD2:AD46").Value = 71
$ws.Range("AE2:AE46").Value = 0
